$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "Auditory Skills Checklist" resource row (row 10). Deleting the
# entire row shifts all subsequent rows up by one and Excel automatically
# compacts the shared-strings table for the now-unused entries.
$ws.Rows.Item(10).Delete()
